# Normalize the "Recorded By" (column G) contributor lists.
#
# Each G-column cell holds a comma-separated list of recorder identities,
# e.g. "dnasr281@gmail.com, System" or "system, backup@backdoor.com, System".
# The canonical ordering keeps a trailing "System" entry but cycles it to the
# front of the list (left-rotation by one) UNLESS the list already starts
# with the backup-account identity "backup@backdoor.com", which is left
# untouched. Single-entry cells, and cells that don't end in "System", are
# also left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ', '

    if ($parts.Count -lt 2) {
        continue
    }

    $first = $parts[0]
    $last = $parts[$parts.Count - 1]

    if ($last -ne 'System') {
        continue
    }

    if ($first -eq 'backup@backdoor.com') {
        continue
    }

    $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
    $newText = $rotated -join ', '

    $cell.Value = $newText
}
